$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the refreshed "Price" values happen to look like plain numbers
# (e.g. "1.00", "7.00"), but every Price/Volume cell in this sheet is stored as
# text. Temporarily switching the cell to a Text number format forces Excel to
# keep the assigned string as text instead of silently parsing it into a
# number; re-applying the built-in "Normal" style afterward drops the cell back
# to the workbook default formatting (no explicit style), matching the original.
function Set-TextValue($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '68.941.00'
$ws.Range("E2").Value = '  -0.24%  '
$ws.Range("D3").Value = '3.766.05'
$ws.Range("E3").Value = '  -1.58%  '
Set-TextValue "D4" '1.00'
$ws.Range("E4").Value = '  +0.69%  '
Set-TextValue "D5" '628.53'
$ws.Range("E5").Value = '  +0.28%  '
Set-TextValue "D6" '165.66'
$ws.Range("E6").Value = '  +0.26%  '
$ws.Range("D7").Value = '3.764.31'
$ws.Range("E7").Value = '  -1.64%  '
$ws.Range("E8").Value = '  +0.03%  '
Set-TextValue "D9" '0.519'
$ws.Range("E9").Value = '  -0.12%  '
$ws.Range("E10").Value = '  -2.14%  '
Set-TextValue "D11" '0.456'
$ws.Range("E11").Value = '  +0.27%  '
$ws.Range("E12").Value = '  +1.26%  '
$ws.Range("E13").Value = '  -4.70%  '
Set-TextValue "D14" '34.73'
$ws.Range("E14").Value = '  -3.75%  '
$ws.Range("D15").Value = '4.398.14'
$ws.Range("E15").Value = '  -1.60%  '
$ws.Range("D16").Value = '3.764.20'
$ws.Range("E16").Value = '  -1.94%  '
$ws.Range("D17").Value = '68.929.28'
$ws.Range("E17").Value = '  -0.25%  '
Set-TextValue "D18" '17.63'
$ws.Range("E18").Value = '  -3.88%  '
$ws.Range("E19").Value = '  -0.14%  '
Set-TextValue "D20" '7.00'
$ws.Range("E20").Value = '  -2.23%  '
Set-TextValue "D21" '461.77'
$ws.Range("E21").Value = '  -1.29%  '
Set-TextValue "D22" '9.49'
$ws.Range("E22").Value = '  -2.32%  '
Set-TextValue "D23" '0.702'
$ws.Range("E23").Value = '  -1.10%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue "D24" '82.03'
$ws.Range("E24").Value = '  -2.21%  '
$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue "D25" '0.0000143'
$ws.Range("E25").Value = '  -6.37%  '
Set-TextValue "D26" '12.08'
$ws.Range("E26").Value = '  +0.00%  '
Set-TextValue "D27" '2.11'
$ws.Range("E27").Value = '  -2.02%  '
Set-TextValue "D28" '10.09'
$ws.Range("E28").Value = '  +0.14%  '
$ws.Range("D30").Value = '3.920.70'
$ws.Range("E30").Value = '  -1.49%  '
Set-TextValue "D31" '2.27'
$ws.Range("E31").Value = '  +0.37%  '
Set-TextValue "D32" '2.66'
$ws.Range("E32").Value = '  -0.17%  '
Set-TextValue "D33" '7.05'
$ws.Range("E33").Value = '  -4.21%  '
Set-TextValue "D34" '28.35'
$ws.Range("E34").Value = '  -2.97%  '
Set-TextValue "D35" '0.174'
$ws.Range("E35").Value = '  +17.01%  '
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("D37").Value = '3.717.78'
$ws.Range("E37").Value = '  -1.46%  '
Set-TextValue "D38" '8.88'
$ws.Range("E38").Value = '  -2.75%  '
$ws.Range("E39").Value = '  -1.75%  '
Set-TextValue "D40" '3.28'
$ws.Range("E40").Value = '  -1.33%  '
Set-TextValue "D41" '5.77'
$ws.Range("E41").Value = '  -2.73%  '
Set-TextValue "D42" '1.00'
$ws.Range("E42").Value = '  +0.08%  '
Set-TextValue "D43" '0.960'
$ws.Range("E43").Value = '  -2.16%  '
$ws.Range("E44").Value = '  -0.02%  '
Set-TextValue "D45" '156.77'
$ws.Range("E45").Value = '  +0.41%  '
$ws.Range("E46").Value = '  +3.76%  '
$ws.Range("E47").Value = '  +0.27%  '
Set-TextValue "D48" '46.93'
$ws.Range("E48").Value = '  +0.31%  '
Set-TextValue "D49" '42.81'
$ws.Range("E49").Value = '  -0.19%  '
Set-TextValue "D50" '0.293'
$ws.Range("E50").Value = '  -2.96%  '
Set-TextValue "D51" '8.33'
$ws.Range("E51").Value = '  -1.52%  '
